# The underlying records in rows 2-11 (minus row 4, which is untouched) were
# re-ordered / re-keyed upstream. This reproduces that by rotating the whole
# row content among the affected rows, in three independent 3-cycles:
#   rows 2 -> 3 -> 5 -> 2
#   rows 6 -> 8 -> 11 -> 6
#   rows 7 -> 9 -> 10 -> 7
# i.e. "new row X" gets the content that used to live in "old row Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row <- source row (content moves from src to dest)
$pairs = @(
    @{dest=2;  src=3},
    @{dest=3;  src=5},
    @{dest=5;  src=2},

    @{dest=6;  src=8},
    @{dest=8;  src=11},
    @{dest=11; src=6},

    @{dest=7;  src=9},
    @{dest=9;  src=10},
    @{dest=10; src=7}
)

# Columns that actually carry record-identifying data (A,B,D,E,F,G,H,Q,R,Y,AA).
# Columns in between (C,I,P,S,T,U,V,W,Z,AB,AD,AE,AG,AT,AW,AX, etc.) are the
# same for every row inside a given cycle, so they don't need touching.
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 25, 27)

# Y (25) and AA (27) hold dates formatted as plain text ("2023-08-29"); a
# bare string assignment would get auto-recognised as a real date serial,
# so those two columns need the text-literal treatment below.
$dateTextCols = @(25, 27)

# Snapshot every source row's values first, since a row can be both a
# source and a destination in the same rotation.
$buffer = @{}
foreach ($p in $pairs) {
    $srcRow = $p.src
    if (-not $buffer.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($c in $cols) {
            $rowVals[$c] = $ws.Cells.Item($srcRow, $c).Value2
        }
        $buffer[$srcRow] = $rowVals
    }
}

# Now write the buffered values into their destination rows.
foreach ($p in $pairs) {
    $destRow = $p.dest
    $srcRow = $p.src
    $rowVals = $buffer[$srcRow]
    foreach ($c in $cols) {
        $val = $rowVals[$c]
        $cell = $ws.Cells.Item($destRow, $c)
        if (($dateTextCols -contains $c) -and ($val -ne $null) -and ($val -ne "")) {
            # Force text (leading apostrophe = "treat as text" like typing
            # it into Excel), then drop the resulting quote-prefix style so
            # the cell format stays plain/default.
            $cell.Value2 = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value2 = $val
        }
    }
}

# Column AF ("Bestämningsmetod") holds a blank placeholder on exactly one
# row per 6/7/8-cycle group; that placeholder moves together with the rest
# of the row's content, so re-home it the same way the rows above were
# rotated (AF8 -> AF6, AF10 -> AF9) and clear it from the rows it leaves.
$ws.Range("AF6").Formula = '=""'
$ws.Range("AF9").Formula = '=""'
$ws.Range("AF8").ClearContents()
$ws.Range("AF10").ClearContents()
